# Auto-generated script: update Moogle Profits market-price derived columns
# per sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match refreshed
# Universalis market data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 440.2
$ws.Range("I11").Value = 440.2
$ws.Range("K11").Value = 440.2
$ws.Range("M11").Value = -300.2
$ws.Range("H17").Value = 2192.8
$ws.Range("J17").Value = 2192.8
$ws.Range("L17").Value = 6578.400000000001
$ws.Range("N17").Value = -6914.400000000001
$ws.Range("H34").Value = 8159.6
$ws.Range("I34").Value = 8159.6
$ws.Range("K34").Value = 8159.6
$ws.Range("M34").Value = -7956.6
$ws.Range("H36").Value = 8159.6
$ws.Range("I36").Value = 8159.6
$ws.Range("K36").Value = 8159.6
$ws.Range("M36").Value = -7444.6
$ws.Range("H38").Value = 4282.64
$ws.Range("I38").Value = 1231.7693
$ws.Range("J38").Value = 7587.75
$ws.Range("K38").Value = 3695.3079
$ws.Range("L38").Value = 22763.25
$ws.Range("M38").Value = -3323.3079
$ws.Range("N38").Value = -23507.25
$ws.Range("H40").Value = 4247.5557
$ws.Range("I40").Value = 1132
$ws.Range("J40").Value = 5805.3335
$ws.Range("K40").Value = 1132
$ws.Range("L40").Value = 5805.3335
$ws.Range("M40").Value = -957
$ws.Range("N40").Value = -6155.3335
$ws.Range("H74").Value = 27281428
$ws.Range("I74").Value = 50007450
$ws.Range("J74").Value = 10199.8
$ws.Range("K74").Value = 50007450
$ws.Range("L74").Value = 10199.8
$ws.Range("M74").Value = -50006514
$ws.Range("N74").Value = -12071.8
$ws.Range("H77").Value = 27281428
$ws.Range("I77").Value = 50007450
$ws.Range("J77").Value = 10199.8
$ws.Range("K77").Value = 250037250
$ws.Range("L77").Value = 50999
$ws.Range("M77").Value = -250032570
$ws.Range("N77").Value = -60359
$ws.Range("H116").Value = 19061.75
$ws.Range("I116").Value = 18999.8
$ws.Range("K116").Value = 18999.8
$ws.Range("M116").Value = -15557.8
$ws.Range("H132").Value = 2243
$ws.Range("I132").Value = 2178.3513
$ws.Range("K132").Value = 6535.053899999999
$ws.Range("M132").Value = -4005.053899999999
$ws.Range("H134").Value = 76208
$ws.Range("J134").Value = 76208
$ws.Range("L134").Value = 76208
$ws.Range("N134").Value = -86348
$ws.Range("H137").Value = 1623.0847
$ws.Range("I137").Value = 1575.6909
$ws.Range("J137").Value = 2274.75
$ws.Range("K137").Value = 4727.072700000001
$ws.Range("L137").Value = 6824.25
$ws.Range("M137").Value = -2177.072700000001
$ws.Range("N137").Value = -11924.25
$ws.Range("H138").Value = 4190.623
$ws.Range("I138").Value = 3439.5715
$ws.Range("J138").Value = 5358.926
$ws.Range("K138").Value = 10318.7145
$ws.Range("L138").Value = 16076.778
$ws.Range("M138").Value = -5178.7145
$ws.Range("N138").Value = -26356.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7007.9575
$ws.Range("I32").Value = 3670.742
$ws.Range("K32").Value = 3670.742
$ws.Range("M32").Value = -3383.742
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("H61").Value = 6471.28
$ws.Range("J61").Value = 10399.6
$ws.Range("L61").Value = 10399.6
$ws.Range("N61").Value = -10823.6
$ws.Range("H122").Value = 2201.7222
$ws.Range("I122").Value = 1999.8541
$ws.Range("K122").Value = 5999.5623
$ws.Range("M122").Value = -3549.5623
$ws.Range("H132").Value = 3641.9268
$ws.Range("I132").Value = 2225.2856
$ws.Range("K132").Value = 6675.8568
$ws.Range("M132").Value = -4145.8568
$ws.Range("H136").Value = 6471.28
$ws.Range("J136").Value = 10399.6
$ws.Range("L136").Value = 31198.8
$ws.Range("N136").Value = -36298.8
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 67415.8
$ws.Range("I82").Value = 20859.6
$ws.Range("J82").Value = 113972
$ws.Range("K82").Value = 20859.6
$ws.Range("L82").Value = 113972
$ws.Range("M82").Value = -20476.6
$ws.Range("N82").Value = -114738
$ws.Range("H85").Value = 67415.8
$ws.Range("I85").Value = 20859.6
$ws.Range("J85").Value = 113972
$ws.Range("K85").Value = 20859.6
$ws.Range("L85").Value = 113972
$ws.Range("M85").Value = -19533.6
$ws.Range("N85").Value = -116624
$ws.Range("H105").Value = 3347.3901
$ws.Range("I105").Value = 3209.4707
$ws.Range("K105").Value = 3209.4707
$ws.Range("M105").Value = -1462.4707
$ws.Range("H107").Value = 3305.3713
$ws.Range("I107").Value = 3393.0293
$ws.Range("K107").Value = 3393.0293
$ws.Range("M107").Value = -1473.0293
$ws.Range("H134").Value = 2413.8076
$ws.Range("I134").Value = 1510.8695
$ws.Range("J134").Value = 9336.333000000001
$ws.Range("K134").Value = 4532.6085
$ws.Range("L134").Value = 28008.999
$ws.Range("M134").Value = -1997.6085
$ws.Range("N134").Value = -33078.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5027.476
$ws.Range("I31").Value = 1959.8889
$ws.Range("J31").Value = 23433
$ws.Range("K31").Value = 1959.8889
$ws.Range("L31").Value = 23433
$ws.Range("M31").Value = -1664.8889
$ws.Range("N31").Value = -24023
$ws.Range("H34").Value = 5027.476
$ws.Range("I34").Value = 1959.8889
$ws.Range("J34").Value = 23433
$ws.Range("K34").Value = 1959.8889
$ws.Range("L34").Value = 23433
$ws.Range("M34").Value = -1757.8889
$ws.Range("N34").Value = -23837
$ws.Range("H99").Value = 2580.158
$ws.Range("I99").Value = 2560.2354
$ws.Range("J99").Value = 2749.5
$ws.Range("K99").Value = 2560.2354
$ws.Range("L99").Value = 2749.5
$ws.Range("M99").Value = -1062.2354
$ws.Range("N99").Value = -5745.5
$ws.Range("H126").Value = 2580.158
$ws.Range("I126").Value = 2560.2354
$ws.Range("J126").Value = 2749.5
$ws.Range("K126").Value = 7680.706200000001
$ws.Range("L126").Value = 8248.5
$ws.Range("M126").Value = -5210.706200000001
$ws.Range("N126").Value = -13188.5
$ws.Range("H132").Value = 3445.2632
$ws.Range("I132").Value = 2352.5925
$ws.Range("J132").Value = 6127.273
$ws.Range("K132").Value = 7057.7775
$ws.Range("L132").Value = 18381.819
$ws.Range("M132").Value = -4527.7775
$ws.Range("N132").Value = -23441.819
$ws.Range("H134").Value = 11630984
$ws.Range("I134").Value = 13890953
$ws.Range("K134").Value = 41672859
$ws.Range("M134").Value = -41670324

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 147.72223
$ws.Range("I38").Value = 104.27273
$ws.Range("J38").Value = 216
$ws.Range("K38").Value = 312.81819
$ws.Range("L38").Value = 648
$ws.Range("M38").Value = 34.18181000000004
$ws.Range("N38").Value = -1342
$ws.Range("H98").Value = 1366.4
$ws.Range("J98").Value = 1343.375
$ws.Range("L98").Value = 4030.125
$ws.Range("N98").Value = -7026.125
$ws.Range("H107").Value = 576
$ws.Range("I107").Value = 288.5
$ws.Range("J107").Value = 633.5
$ws.Range("K107").Value = 865.5
$ws.Range("L107").Value = 1900.5
$ws.Range("M107").Value = 1054.5
$ws.Range("N107").Value = -5740.5
$ws.Range("H132").Value = 3446.375
$ws.Range("I132").Value = 3356.3333
$ws.Range("J132").Value = 3500.4
$ws.Range("K132").Value = 30206.9997
$ws.Range("L132").Value = 31503.6
$ws.Range("M132").Value = -27676.9997
$ws.Range("N132").Value = -36563.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10000
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 10000
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H80").Value = 4987.6313
$ws.Range("I80").Value = 3231.25
$ws.Range("K80").Value = 3231.25
$ws.Range("M80").Value = -2233.25
$ws.Range("H83").Value = 4987.6313
$ws.Range("I83").Value = 3231.25
$ws.Range("K83").Value = 16156.25
$ws.Range("M83").Value = -11164.25
$ws.Range("H110").Value = 56702
$ws.Range("J110").Value = 56702
$ws.Range("L110").Value = 56702
$ws.Range("N110").Value = -64882
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 43886.92
$ws.Range("I7").Value = 45611.418
$ws.Range("K7").Value = 45611.418
$ws.Range("M7").Value = -45499.418
$ws.Range("H55").Value = 761.36365
$ws.Range("I55").Value = 100.75
$ws.Range("J55").Value = 1554.1
$ws.Range("K55").Value = 100.75
$ws.Range("L55").Value = 1554.1
$ws.Range("M55").Value = 72.25
$ws.Range("N55").Value = -1900.1
$ws.Range("H98").Value = 100000
$ws.Range("J98").Value = 100000
$ws.Range("L98").Value = 100000
$ws.Range("N98").Value = -105990
$ws.Range("H126").Value = 43886.92
$ws.Range("I126").Value = 45611.418
$ws.Range("K126").Value = 136834.254
$ws.Range("M126").Value = -134364.254
$ws.Range("H132").Value = 3320
$ws.Range("I132").Value = 1789.6923
$ws.Range("J132").Value = 6635.6665
$ws.Range("K132").Value = 5369.0769
$ws.Range("L132").Value = 19906.9995
$ws.Range("M132").Value = -2839.0769
$ws.Range("N132").Value = -24966.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2913.2903
$ws.Range("I132").Value = 2161.9614
$ws.Range("K132").Value = 6485.8842
$ws.Range("M132").Value = -3955.8842
$ws.Range("H136").Value = 4982.5293
$ws.Range("I136").Value = 4970.4
$ws.Range("K136").Value = 14911.2
$ws.Range("M136").Value = -12361.2
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 123665.664
$ws.Range("J138").Value = 123665.664
$ws.Range("L138").Value = 123665.664
$ws.Range("N138").Value = -133945.664
$ws.Range("N137").ClearContents()
